# Generate Report for Handoff
# Adds two new files (6343ef63-df64-431c-967b-76b16363b5c6.md and
# 88b0901b-12f9-4eb2-a346-43a61ba2770a.md) as new rows to the Overview,
# zh-cn and de-de report tables.

$wb = $excel.ActiveWorkbook

$file1 = "6343ef63-df64-431c-967b-76b16363b5c6.md"
$file2 = "88b0901b-12f9-4eb2-a346-43a61ba2770a.md"

$hash1 = "a62f68a5e108c99348e9a9edb6c63101cd737acb"
$hash2 = "8cc47e9662fac810c7bd2df1c6c67075b9836601"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): columns A..G
#   A File Name, B Path And Name, C Extension, D Publish URL,
#   E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)

$tblOverview.ListRows.Add() | Out-Null
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $file1
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$hash1/e2e/$file1", [System.Type]::Missing, [System.Type]::Missing, "e2e\" + $file1)
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-29 08:43:43"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = $file2
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$hash2/e2e/$file2", [System.Type]::Missing, [System.Type]::Missing, "e2e\" + $file2)
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-29 08:43:43"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) and "de-de" (sheet3): columns A..P
#   A Source File Name, B File Extension, C Status, D Source Path,
#   E Priority, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, I Latest Target File,
#   J Latest Handback File, K Latest Handback DateTime,
#   L Reference Tokens, M To be localized, N Dependency From,
#   O Has metadata, P Error Detail
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; H4 = "2016-08-29 08:43:39"; H5 = "2016-08-29 08:43:39" },
    @{ Name = "de-de"; H4 = "2016-08-29 08:43:43"; H5 = "2016-08-29 08:43:43" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)
    $tbl = $ws.ListObjects.Item(1)

    $tbl.ListRows.Add() | Out-Null
    $tbl.ListRows.Add() | Out-Null

    # --- row 4: file1 ---
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$hash1/e2e/$file1", [System.Type]::Missing, [System.Type]::Missing, $file1)
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = "e2e"
    $ws.Range("E4").Value = "ht"
    $ws.Range("F4").Value = "False"
    $ws.Range("G4").Value = "$file1.$hash1.$($loc.Name).xlf"
    $ws.Range("H4").Value = $loc.H4
    $ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("K4").Value = "0001-01-01 00:00:00"
    $ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("M4").Value = "True"
    $ws.Range("O4").Value = "False"

    # --- row 5: file2 ---
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$hash2/e2e/$file2", [System.Type]::Missing, [System.Type]::Missing, $file2)
    $ws.Range("B5").Value = ".md"
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Range("D5").Value = "e2e"
    $ws.Range("E5").Value = "ht"
    $ws.Range("F5").Value = "False"
    $ws.Range("G5").Value = "$file2.$hash2.$($loc.Name).xlf"
    $ws.Range("H5").Value = $loc.H5
    $ws.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("K5").Value = "0001-01-01 00:00:00"
    $ws.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("M5").Value = "True"
    $ws.Range("O5").Value = "False"
}
